# SwaadSutra_Daily_2026-01-13.xlsx update
# A brand new order (Order ID 6, "Wheat Chapati x1" placed by Pooja) came in
# at 2026-01-13 16:41 and is inserted at the top of the "Daily Orders" log,
# pushing every existing order down by one row. The "Summary" and
# "Items Breakdown" sheets are refreshed to reflect the new totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily Orders - insert the new order as the new row 2
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")

# Shift all existing data rows (2..6) down by one to make room for the
# newest order at the top of the list.
$orders.Rows.Item(2).Insert()

$orders.Range("A2").Value = 6
$orders.Range("B2").Value = "2026-01-13 16:41"
$orders.Range("C2").Value = "Pooja"
$orders.Range("D2").Value = "saf"

# Phone numbers / dates that look numeric must be kept as plain text,
# exactly like the rest of the sheet.
$orders.Range("E2").NumberFormat = "@"
$orders.Range("E2").Value = "9096648553"

$orders.Range("F2").Value = "Wheat Chapati x1"
$orders.Range("G2").Value = 15
$orders.Range("H2").Value = "NEW"
$orders.Range("I2").Value = "PENDING"

$orders.Range("J2").NumberFormat = "@"
$orders.Range("J2").Value = "2026-01-15"

$orders.Range("K2").Value = "10:12"

# ---------------------------------------------------------------------
# 2) Summary - recalculated totals for the day
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A2").Value = 6    # Total Orders
$summary.Range("B2").Value = 5    # New
$summary.Range("C2").Value = 1    # Cooking
$summary.Range("D2").Value = 0    # Ready
$summary.Range("E2").Value = 0    # Delivered
$summary.Range("F2").Value = 0    # Cancelled
$summary.Range("G2").Value = 180  # Total Revenue
$summary.Range("H2").Value = 0    # Paid Amount

# ---------------------------------------------------------------------
# 3) Items Breakdown - recomputed per-item quantity / revenue, ordered by
#    the most recent order each item appeared in.
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")

$items.Range("A2").Value = "Wheat Chapati"
$items.Range("B2").Value = 2
$items.Range("C2").Value = 30

$items.Range("A3").Value = "Til Poli"
$items.Range("B3").Value = 2
$items.Range("C3").Value = 60

$items.Range("A4").Value = "Onion Pakoda (Kanda Bhaje)"
$items.Range("B4").Value = 1
$items.Range("C4").Value = 60

$items.Range("A5").Value = "Pohe"
$items.Range("B5").Value = 1
$items.Range("C5").Value = 30
